# DeveloperGuide: remove all references to `BrowserPanel`
#
# The `BrowserPanel` UI class no longer exists, so its box is removed from
# the UiClassDiagram along with the two connector lines that were glued to
# it. The boxes/connectors below the removed BrowserPanel box shift up to
# close the gap, and a couple of surrounding shapes are resized/repositioned
# to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$EMU = 12700.0

# Shape.Left/Top/Width/Height are (single-precision) points; the host
# truncates the EMU conversion rather than rounding, so nudge each value up
# by half an EMU before dividing to land on the exact target EMU.
function EmuToPt($emu) {
    return ([double]$emu + 0.5) / $EMU
}

function Set-ShapeFrame($shapes, $id, $x, $y, $cx, $cy) {
    $sh = Get-ShapeById $shapes $id
    $sh.Left = EmuToPt $x
    $sh.Top = EmuToPt $y
    $sh.Width = EmuToPt $cx
    $sh.Height = EmuToPt $cy
}

$shapes = $s.Shapes

# --- Remove the BrowserPanel box (id 34) and its two glued connectors ---
(Get-ShapeById $shapes 34).Delete()     # "BrowserPanel" rectangle
(Get-ShapeById $shapes 44).Delete()     # connector: Flowchart Decision 39 -> BrowserPanel
(Get-ShapeById $shapes 82).Delete()     # connector: Isosceles Triangle 16 -> BrowserPanel

# --- Resize the outer "UI" container now that BrowserPanel's row is gone ---
Set-ShapeFrame $shapes 23 1217465 1447800 4917083 3581393

# --- Reposition/resize the "Logic" bar label alongside it ---
Set-ShapeFrame $shapes 22 5897465 2271241 1974930 328045

# --- Shift the remaining stacked boxes up to close the gap left behind ---
Set-ShapeFrame $shapes 35 2592527 4182760 1093635 236841   # StatusBarFooter
Set-ShapeFrame $shapes 36 2592526 3610961 1093635 236841   # PersonListPanel
Set-ShapeFrame $shapes 37 3839323 3847802 1040906 236841   # PersonCard
Set-ShapeFrame $shapes 38 2592528 4585001 1093635 236841   # HelpWindow

# --- Re-route the connectors attached to those boxes to their new spots ---
Set-ShapeFrame $shapes 47 2073648 3210503 861357 176400
Set-ShapeFrame $shapes 50 1787748 3496402 1433156 176401
Set-ShapeFrame $shapes 53 1374846 3485740 2018094 417270
Set-ShapeFrame $shapes 77 4364988 2801241 1680223 649740
Set-ShapeFrame $shapes 91 3597466 2371782 2018095 1840702
Set-ShapeFrame $shapes 94 3399355 2572808 2417422 1843806
Set-ShapeFrame $shapes 140 3886374 2085787 1443382 1843808
Set-ShapeFrame $shapes 137 3430123 3557022 118421 699979

# --- Other nearby shapes that moved up with the stack ---
Set-ShapeFrame $shapes 117 6213739 4179377 1371599 328045  # "Model" bar label
Set-ShapeFrame $shapes 144 5431573 4107138 229325 160062
Set-ShapeFrame $shapes 118 4114799 4091709 2642195 101600
